$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The edit reshuffles the species-observation data held in rows 18-21:
#   - Row 18 receives the data that used to live in row 21 (Ullticka)
#   - Row 19 receives the data that used to live in row 20 (Lunglav)
#   - Row 20 receives the data that used to live in row 18 (Tretåig hackspett)
#   - Row 21 receives the data that used to live in row 19 (Tretåig hackspett)
# Columns D,I,P,S,T,U,V,W,Y,AA,AD,AE,AG,AT,AW,AX,AY are identical across the
# four rows, so they are left untouched. Only A,B,E,F,G,H,Q,R,Z,AB change
# value, and K,L,M,N,AC (age/sex/activity/method/public-comment, which only
# apply to the woodpecker sighting) move from rows 18/19 to rows 20/21.
# ---------------------------------------------------------------------------

# Row 18 <- old row 21 (Ullticka / Phellinidium ferrugineofuscum)
$ws.Range("A18").Value = 130979947
$ws.Range("B18").Value = 91808
$ws.Range("E18").Value = 1202
$ws.Range("F18").Value = "Ullticka"
$ws.Range("G18").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H18").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("Q18").Value = 590591
$ws.Range("R18").Value = 6963354
$ws.Range("Z18").Value = "09:45"
$ws.Range("AB18").Value = "09:45"
$ws.Range("K18").ClearContents()
$ws.Range("L18").ClearContents()
$ws.Range("M18").ClearContents()
$ws.Range("N18").ClearContents()
$ws.Range("AC18").ClearContents()

# Row 19 <- old row 20 (Lunglav / Lobaria pulmonaria)
$ws.Range("A19").Value = 130979897
$ws.Range("B19").Value = 80348
$ws.Range("E19").Value = 6458
$ws.Range("F19").Value = "Lunglav"
$ws.Range("G19").Value = "Lobaria pulmonaria"
$ws.Range("H19").Value = "(L.) Hoffm."
$ws.Range("Q19").Value = 590726
$ws.Range("R19").Value = 6963153
$ws.Range("Z19").Value = "13:24"
$ws.Range("AB19").Value = "13:24"
$ws.Range("K19").ClearContents()
$ws.Range("L19").ClearContents()
$ws.Range("M19").ClearContents()
$ws.Range("N19").ClearContents()
$ws.Range("AC19").ClearContents()

# Row 20 <- old row 18 (Tretåig hackspett / Picoides tridactylus)
$ws.Range("A20").Value = 130979946
$ws.Range("B20").Value = 57884
$ws.Range("E20").Value = 100109
$ws.Range("F20").Value = "Tretåig hackspett"
$ws.Range("G20").Value = "Picoides tridactylus"
$ws.Range("H20").Value = "(Linnaeus, 1758)"
$ws.Range("Q20").Value = 590605
$ws.Range("R20").Value = 6963364
$ws.Range("Z20").Value = "09:47"
$ws.Range("AB20").Value = "09:47"
$ws.Range("K20").Value = ""
$ws.Range("L20").Value = ""
$ws.Range("M20").Value = "färska spår"
$ws.Range("N20").Value = ""
$ws.Range("AC20").Value = "färska ringhack på gran"

# Row 21 <- old row 19 (Tretåig hackspett / Picoides tridactylus)
$ws.Range("A21").Value = 130979899
$ws.Range("B21").Value = 57884
$ws.Range("E21").Value = 100109
$ws.Range("F21").Value = "Tretåig hackspett"
$ws.Range("G21").Value = "Picoides tridactylus"
$ws.Range("H21").Value = "(Linnaeus, 1758)"
$ws.Range("Q21").Value = 590850
$ws.Range("R21").Value = 6963133
$ws.Range("Z21").Value = "13:16"
$ws.Range("AB21").Value = "13:16"
$ws.Range("K21").Value = ""
$ws.Range("L21").Value = ""
$ws.Range("M21").Value = "färska spår"
$ws.Range("N21").Value = ""
$ws.Range("AC21").Value = "färska ringhack på gran"

Write-Host "Rows 18-21 updated"
